# Insert a new row at position 968 (shifts existing rows 968:1036 down to 969:1037)
# and populate it with a new price observation record.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(968).Insert()

$ws.Range("A968").Value = 10
$ws.Range("B968").Value = 'Vega Modelo de Temuco'
$ws.Range("C968").Value = 'La Araucanía'
$ws.Range("D968").Value = 45265
$ws.Range("E968").Value = 9
$ws.Range("F968").Value = 100112032
$ws.Range("G968").Value = 'Zapallo italiano'
$ws.Range("H968").Value = 'Bola 8'
$ws.Range("I968").Value = 'Primera'
$ws.Range("J968").Value = 50
$ws.Range("K968").Value = 15000
$ws.Range("L968").Value = 15000
$ws.Range("M968").Value = 15000
$ws.Range("N968").Value = '$/caja 50 unidades'
$ws.Range("O968").Value = "Región de O'Higgins"
$ws.Range("P968").Value = 300
$ws.Range("Q968").Value = 50
$ws.Range("R968").Value = 'Hortaliza'
